# "Actualizo pib y VA al tercer trimestre de 2023"
#
# 1) VAB sheet: refresh B2:C18 ("VA" contributions and variation %).
# 2) Producto sheet: refresh column C (PIB DESEST) for rows 2:79 and
#    append the new 2023-III row (row 80).
# 3) Reflect the new active tab / selections as left by the author.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) VAB sheet — updated values for B2:C18
# ---------------------------------------------------------------------
$wsVab = $wb.Worksheets.Item("VAB")

$vabData = @(
    @(2, 601421.19613162649, -0.43489291804220676),
    @(3, 31595.916082538191, -7.6343099569625039),
    @(4, 3858.3631153696829, 2.9435981467078154),
    @(5, 26267.492834246732, 5.871793690677185),
    @(6, 122612.92120666231, -3.6795710755112565),
    @(7, 13209.273653523698, 2.6041227585157145),
    @(8, 23620.858649502094, -0.12885601193497465),
    @(9, 97940.134100454714, -0.51240500487524931),
    @(10, 11039.319275525062, 7.2162068777362309),
    @(11, 54289.274559934063, -0.15009459778109191),
    @(12, 24835.947824718762, -0.96303111987976475),
    @(13, 82209.628988706929, 1.0399730239221583),
    @(14, 35867.316184027521, 2.4423649945773063),
    @(15, 29487.346840138496, 2.607571996301461),
    @(16, 24558.162985526862, 1.7285235675380983),
    @(17, 16023.739564855718, 1.0597913614180721),
    @(18, 4005.5002658956105, -3.7218612669402407)
)

foreach ($row in $vabData) {
    $r = $row[0]
    $wsVab.Cells.Item($r, 2).Value = $row[1]
    $wsVab.Cells.Item($r, 3).Value = $row[2]
}

# ---------------------------------------------------------------------
# 2) Producto sheet — updated column C (PIB DESEST) for rows 2:79
# ---------------------------------------------------------------------
$wsProducto = $wb.Worksheets.Item("Producto")

$productoC = @(
    @(2, 475908.01430203719),
    @(3, 470258.31897942215),
    @(4, 493778.21076235082),
    @(5, 500516.23684607708),
    @(6, 514850.78048926021),
    @(7, 527217.04692234145),
    @(8, 530089.29816912254),
    @(9, 540066.64443293249),
    @(10, 554209.42435356788),
    @(11, 561322.81502015446),
    @(12, 577233.71960137994),
    @(13, 589431.65790782252),
    @(14, 603133.40302888909),
    @(15, 616014.20510520146),
    @(16, 624860.00809330563),
    @(17, 643762.3943569666),
    @(18, 649653.75390456943),
    @(19, 653109.86116908432),
    @(20, 658442.76983749284),
    @(21, 627498.25405370165),
    @(22, 604368.91495259106),
    @(23, 590301.51436368749),
    @(24, 614722.58336871513),
    @(25, 626098.492966505),
    @(26, 644922.61968795315),
    @(27, 673884.146169079),
    @(28, 677943.34671504016),
    @(29, 685344.6051951237),
    @(30, 700946.23223711469),
    @(31, 710462.4782490026),
    @(32, 716501.86977298337),
    @(33, 715215.80862330692),
    @(34, 707585.26657786965),
    @(35, 684114.06975628051),
    @(36, 705462.50859424984),
    @(37, 716782.11290739593),
    @(38, 715570.03785572597),
    @(39, 721318.15144984378),
    @(40, 725469.43745116564),
    @(41, 719270.79445452429),
    @(42, 708041.5307386308),
    @(43, 703312.6787777018),
    @(44, 697146.70795409917),
    @(45, 700723.26638302929),
    @(46, 710648.05387179856),
    @(47, 729189.84750436153),
    @(48, 727321.93472757621),
    @(49, 718788.75044842588),
    @(50, 713642.99007153255),
    @(51, 701335.42587019689),
    @(52, 703166.588613534),
    @(53, 707766.38983538363),
    @(54, 714972.09359957289),
    @(55, 721394.01533615473),
    @(56, 731146.34352879832),
    @(57, 738047.33858676918),
    @(58, 734287.63460992684),
    @(59, 702979.52002783224),
    @(60, 701760.01104200038),
    @(61, 690482.60623829148),
    @(62, 693076.13247394061),
    @(63, 696153.17451175139),
    @(64, 697087.14028623956),
    @(65, 686578.76707990898),
    @(66, 655959.36861962499),
    @(67, 563046.74878187256),
    @(68, 625286.7337214232),
    @(69, 654072.29347866599),
    @(70, 677896.51590083656),
    @(71, 678957.75254003098),
    @(72, 696419.82516886073),
    @(73, 712866.07543587033),
    @(74, 718443.35816139122),
    @(75, 730649.11533843062),
    @(76, 733274.29089024279),
    @(77, 720873.5485705015),
    @(78, 726708.20478132647),
    @(79, 706885.32725060661)
)

foreach ($row in $productoC) {
    $r = $row[0]
    $wsProducto.Cells.Item($r, 3).Value = $row[1]
}

# New row 80: 2023 - III quarter
$wsProducto.Cells.Item(80, 1).Value = 2023
$wsProducto.Cells.Item(80, 2).Value = "III"
$wsProducto.Cells.Item(80, 3).Value = 726138.66856733966
$wsProducto.Cells.Item(80, 4).Value = 721353.07659723423

# Reflect the new selection left on the Producto sheet by the author
$wsProducto.Activate() | Out-Null
$wsProducto.Range("D78:D80").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) Leave the VAB sheet as the active tab (matches activeTab change
#    from Aperturas -> VAB) and restore its former selection.
# ---------------------------------------------------------------------
$wsVab.Activate() | Out-Null
$wsVab.Range("C2:C18").Select() | Out-Null
